$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.062088900924398
$ws.Cells.Item(2, 4).Value = 1.066775190774129
$ws.Cells.Item(2, 5).Value = 1.05761722965166
$ws.Cells.Item(2, 6).Value = 1.075572017937336
$ws.Cells.Item(2, 9).Value = 1.05282507205158
$ws.Cells.Item(2, 10).Value = 1.067061224130508
$ws.Cells.Item(2, 11).Value = 1.069484833642417
$ws.Cells.Item(2, 12).Value = 1.060351709029095
$ws.Cells.Item(2, 13).Value = 1.078258234190974
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.06331774159188
$ws.Cells.Item(3, 4).Value = 1.067772924698103
$ws.Cells.Item(3, 5).Value = 1.05867543085846
$ws.Cells.Item(3, 6).Value = 1.076701715084911
$ws.Cells.Item(3, 9).Value = 1.053222933430653
$ws.Cells.Item(3, 10).Value = 1.067943031481684
$ws.Cells.Item(3, 11).Value = 1.070297771097248
$ws.Cells.Item(3, 12).Value = 1.06122316132055
$ws.Cells.Item(3, 13).Value = 1.079204506370166
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.064112700855984
$ws.Cells.Item(4, 4).Value = 1.06841832831858
$ws.Cells.Item(4, 5).Value = 1.059360216827389
$ws.Cells.Item(4, 6).Value = 1.077432804832294
$ws.Cells.Item(4, 9).Value = 1.053479079150611
$ws.Cells.Item(4, 10).Value = 1.068512892741861
$ws.Cells.Item(4, 11).Value = 1.070822983240495
$ws.Cells.Item(4, 12).Value = 1.06178650676341
$ws.Cells.Item(4, 13).Value = 1.079816316361867
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.064446859835625
$ws.Cells.Item(5, 4).Value = 1.068689609687075
$ws.Cells.Item(5, 5).Value = 1.059648115906481
$ws.Cells.Item(5, 6).Value = 1.077740179960075
$ws.Cells.Item(5, 9).Value = 1.053586452756275
$ws.Cells.Item(5, 10).Value = 1.068752289398869
$ws.Cells.Item(5, 11).Value = 1.071043588724869
$ws.Cells.Item(5, 12).Value = 1.062023208226536
$ws.Cells.Item(5, 13).Value = 1.080073404597423
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.064502964196398
$ws.Cells.Item(6, 4).Value = 1.068735156330973
$ws.Cells.Item(6, 5).Value = 1.059696456338601
$ws.Cells.Item(6, 6).Value = 1.07779179109214
$ws.Cells.Item(6, 9).Value = 1.05360446308981
$ws.Cells.Item(6, 10).Value = 1.068792475035522
$ws.Cells.Item(6, 11).Value = 1.071080618023098
$ws.Cells.Item(6, 12).Value = 1.062062943907718
$ws.Cells.Item(6, 13).Value = 1.080116564028674
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.064117166066895
$ws.Cells.Item(7, 4).Value = 1.068421953374439
$ws.Cells.Item(7, 5).Value = 1.059364063688649
$ws.Cells.Item(7, 6).Value = 1.07743691189429
$ws.Cells.Item(7, 9).Value = 1.053480515099667
$ws.Cells.Item(7, 10).Value = 1.068516092250242
$ws.Cells.Item(7, 11).Value = 1.070825931741795
$ws.Cells.Item(7, 12).Value = 1.061789670085543
$ws.Cells.Item(7, 13).Value = 1.079819752044753
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.062504232361557
$ws.Cells.Item(8, 4).Value = 1.067112420907355
$ws.Cells.Item(8, 5).Value = 1.057974842039846
$ws.Cells.Item(8, 6).Value = 1.075953784255589
$ws.Cells.Item(8, 9).Value = 1.052959800059339
$ws.Cells.Item(8, 10).Value = 1.067359386224278
$ws.Cells.Item(8, 11).Value = 1.06975973847258
$ws.Cells.Item(8, 12).Value = 1.060646333276851
$ws.Cells.Item(8, 13).Value = 1.078578133622755
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.059660552935182
$ws.Cells.Item(9, 4).Value = 1.064803311159054
$ws.Cells.Item(9, 5).Value = 1.055527266182552
$ws.Cells.Item(9, 6).Value = 1.073341044753619
$ws.Cells.Item(9, 9).Value = 1.052032276235904
$ws.Cells.Item(9, 10).Value = 1.065315506977774
$ws.Cells.Item(9, 11).Value = 1.067874705571736
$ws.Cells.Item(9, 12).Value = 1.058627432233179
$ws.Cells.Item(9, 13).Value = 1.076386436853937
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.057763649128562
$ws.Cells.Item(10, 4).Value = 1.063262809491461
$ws.Cells.Item(10, 5).Value = 1.053895763530794
$ws.Cells.Item(10, 6).Value = 1.07159963585963
$ws.Cells.Item(10, 9).Value = 1.051407199565406
$ws.Cells.Item(10, 10).Value = 1.063949077355302
$ws.Cells.Item(10, 11).Value = 1.066613747973537
$ws.Cells.Item(10, 12).Value = 1.057278617955635
$ws.Cells.Item(10, 13).Value = 1.074922683793628
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.056941976444258
$ws.Cells.Item(11, 4).Value = 1.062595482308249
$ws.Cells.Item(11, 5).Value = 1.053189341329666
$ws.Cells.Item(11, 6).Value = 1.070845669778339
$ws.Cells.Item(11, 9).Value = 1.051134931101969
$ws.Cells.Item(11, 10).Value = 1.063356470803053
$ws.Cells.Item(11, 11).Value = 1.066066713403869
$ws.Cells.Item(11, 12).Value = 1.056693869389604
$ws.Cells.Item(11, 13).Value = 1.074288225982125
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.056636722978921
$ws.Cells.Item(12, 4).Value = 1.062347564005917
$ws.Cells.Item(12, 5).Value = 1.052926947717614
$ws.Cells.Item(12, 6).Value = 1.07056562307757
$ws.Cells.Item(12, 9).Value = 1.051033556436814
$ws.Cells.Item(12, 10).Value = 1.063136208636865
$ws.Cells.Item(12, 11).Value = 1.065863364167183
$ws.Cells.Item(12, 12).Value = 1.056476560812768
$ws.Cells.Item(12, 13).Value = 1.074052462121572
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.05670220307889
$ws.Cells.Item(13, 4).Value = 1.062400745301493
$ws.Cells.Item(13, 5).Value = 1.052983231921832
$ws.Cells.Item(13, 6).Value = 1.070625693622359
$ws.Cells.Item(13, 9).Value = 1.05105531261095
$ws.Cells.Item(13, 10).Value = 1.063183462055185
$ws.Cells.Item(13, 11).Value = 1.065906990359118
$ws.Cells.Item(13, 12).Value = 1.056523179099483
$ws.Cells.Item(13, 13).Value = 1.074103038725742
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.056916745083726
$ws.Cells.Item(14, 4).Value = 1.062574990193982
$ws.Cells.Item(14, 5).Value = 1.053167651737459
$ws.Cells.Item(14, 6).Value = 1.07082252083779
$ws.Cells.Item(14, 9).Value = 1.051126556382192
$ws.Cells.Item(14, 10).Value = 1.063338266762903
$ws.Cells.Item(14, 11).Value = 1.066049907687052
$ws.Cells.Item(14, 12).Value = 1.056675908783549
$ws.Cells.Item(14, 13).Value = 1.074268739674535
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.057048925106848
$ws.Cells.Item(15, 4).Value = 1.062682342451777
$ws.Cells.Item(15, 5).Value = 1.053281279236026
$ws.Cells.Item(15, 6).Value = 1.070943793845609
$ws.Cells.Item(15, 9).Value = 1.051170419927139
$ws.Cells.Item(15, 10).Value = 1.063433628231343
$ws.Cells.Item(15, 11).Value = 1.066137943028707
$ws.Cells.Item(15, 12).Value = 1.056769996375488
$ws.Cells.Item(15, 13).Value = 1.074370820473709
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.057818174291483
$ws.Cells.Item(16, 4).Value = 1.063307091798304
$ws.Cells.Item(16, 5).Value = 1.053942646946764
$ws.Cells.Item(16, 6).Value = 1.071649675570238
$ws.Cells.Item(16, 9).Value = 1.051425235224744
$ws.Cells.Item(16, 10).Value = 1.063988386908154
$ws.Cells.Item(16, 11).Value = 1.066650031014591
$ws.Cells.Item(16, 12).Value = 1.057317410868673
$ws.Cells.Item(16, 13).Value = 1.0749647770037
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.058300621714148
$ws.Cells.Item(17, 4).Value = 1.063698904775461
$ws.Cells.Item(17, 5).Value = 1.054357512312995
$ws.Cells.Item(17, 6).Value = 1.072092475696356
$ws.Cells.Item(17, 9).Value = 1.05158464360345
$ws.Cells.Item(17, 10).Value = 1.064336121414022
$ws.Cells.Item(17, 11).Value = 1.066970973203602
$ws.Cells.Item(17, 12).Value = 1.057660600369904
$ws.Cells.Item(17, 13).Value = 1.075337177109548
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.0585819962757
$ws.Cells.Item(18, 4).Value = 1.06392741574487
$ws.Cells.Item(18, 5).Value = 1.054599499172869
$ws.Cells.Item(18, 6).Value = 1.072350760899677
$ws.Cells.Item(18, 9).Value = 1.051677468894496
$ws.Cells.Item(18, 10).Value = 1.064538859097683
$ws.Cells.Item(18, 11).Value = 1.067158073903058
$ws.Cells.Item(18, 12).Value = 1.057860709250883
$ws.Cells.Item(18, 13).Value = 1.07555432983664
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.058677932908664
$ws.Cells.Item(19, 4).Value = 1.064005327552497
$ws.Cells.Item(19, 5).Value = 1.054682010972364
$ws.Cells.Item(19, 6).Value = 1.072438830775092
$ws.Cells.Item(19, 9).Value = 1.051709093664788
$ws.Cells.Item(19, 10).Value = 1.064607972191706
$ws.Cells.Item(19, 11).Value = 1.067221853575578
$ws.Cells.Item(19, 12).Value = 1.05792892975532
$ws.Cells.Item(19, 13).Value = 1.075628362837163
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.058248862675351
$ws.Cells.Item(20, 4).Value = 1.063656869744519
$ws.Cells.Item(20, 5).Value = 1.054313000895041
$ws.Cells.Item(20, 6).Value = 1.072044966681572
$ws.Cells.Item(20, 9).Value = 1.051567556622849
$ws.Cells.Item(20, 10).Value = 1.064298822109969
$ws.Cells.Item(20, 11).Value = 1.066936549440509
$ws.Cells.Item(20, 12).Value = 1.057623786402121
$ws.Cells.Item(20, 13).Value = 1.075297228521549
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.056853569145566
$ws.Cells.Item(21, 4).Value = 1.06252368062061
$ws.Cells.Item(21, 5).Value = 1.053113344624988
$ws.Cells.Item(21, 6).Value = 1.070764559867899
$ws.Cells.Item(21, 9).Value = 1.051105583554489
$ws.Cells.Item(21, 10).Value = 1.06329268455841
$ws.Cells.Item(21, 11).Value = 1.066007826419186
$ws.Cells.Item(21, 12).Value = 1.056630936660291
$ws.Cells.Item(21, 13).Value = 1.074219947590813
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.055976015847032
$ws.Cells.Item(22, 4).Value = 1.061810948658118
$ws.Cells.Item(22, 5).Value = 1.052359090074913
$ws.Cells.Item(22, 6).Value = 1.069959573674146
$ws.Cells.Item(22, 9).Value = 1.050813722048373
$ws.Cells.Item(22, 10).Value = 1.062659265350228
$ws.Cells.Item(22, 11).Value = 1.065422996811269
$ws.Cells.Item(22, 12).Value = 1.056006073352793
$ws.Cells.Item(22, 13).Value = 1.073542050751624
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.056441250518677
$ws.Cells.Item(23, 4).Value = 1.062188805561643
$ws.Cells.Item(23, 5).Value = 1.052758933524499
$ws.Cells.Item(23, 6).Value = 1.070386307231641
$ws.Cells.Item(23, 9).Value = 1.050968576369558
$ws.Cells.Item(23, 10).Value = 1.062995131127883
$ws.Cells.Item(23, 11).Value = 1.065733112215189
$ws.Cells.Item(23, 12).Value = 1.056337384347201
$ws.Cells.Item(23, 13).Value = 1.073901470947332
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.058272250446238
$ws.Cells.Item(24, 4).Value = 1.063675863648893
$ws.Cells.Item(24, 5).Value = 1.054333113679963
$ws.Cells.Item(24, 6).Value = 1.072066433937662
$ws.Cells.Item(24, 9).Value = 1.051575277972602
$ws.Cells.Item(24, 10).Value = 1.064315676338004
$ws.Cells.Item(24, 11).Value = 1.066952104364658
$ws.Cells.Item(24, 12).Value = 1.057640421258939
$ws.Cells.Item(24, 13).Value = 1.07531527976042
$ws.Cells.Item(25, 2).Value = 1.019999999999999
$ws.Cells.Item(25, 3).Value = 1.060395900423877
$ws.Cells.Item(25, 4).Value = 1.065400460322701
$ws.Cells.Item(25, 5).Value = 1.056159980889577
$ws.Cells.Item(25, 6).Value = 1.074016422024656
$ws.Cells.Item(25, 9).Value = 1.052273246741164
$ws.Cells.Item(25, 10).Value = 1.065844571201166
$ws.Cells.Item(25, 11).Value = 1.068362779988752
$ws.Cells.Item(25, 12).Value = 1.059149869844058
$ws.Cells.Item(25, 13).Value = 1.076953500528012
